$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 2046554.5   # H17: 1391473.2 -> 2046554.5
$ws.Cells.Item(17, 10).Value = 2046554.5   # J17: 1391473.2 -> 2046554.5
$ws.Cells.Item(17, 12).Value = 6139663.5   # L17: 4174419.6 -> 6139663.5
$ws.Cells.Item(17, 14).Value = -6139999.5   # N17: -4174755.6 -> -6139999.5
$ws.Cells.Item(86, 8).Value = 20836342   # H86: 29415396 -> 20836342
$ws.Cells.Item(86, 9).Value = 33335800   # I86: 55558770 -> 33335800
$ws.Cells.Item(86, 10).Value = 3911.6667   # J86: 4100 -> 3911.6667
$ws.Cells.Item(86, 11).Value = 33335800   # K86: 55558770 -> 33335800
$ws.Cells.Item(86, 12).Value = 3911.6667   # L86: 4100 -> 3911.6667
$ws.Cells.Item(86, 13).Value = -33334677   # M86: -55557647 -> -33334677
$ws.Cells.Item(86, 14).Value = -6157.6667   # N86: -6346 -> -6157.6667
$ws.Cells.Item(89, 8).Value = 20836342   # H89: 29415396 -> 20836342
$ws.Cells.Item(89, 9).Value = 33335800   # I89: 55558770 -> 33335800
$ws.Cells.Item(89, 10).Value = 3911.6667   # J89: 4100 -> 3911.6667
$ws.Cells.Item(89, 11).Value = 166679000   # K89: 277793850 -> 166679000
$ws.Cells.Item(89, 12).Value = 19558.3335   # L89: 20500 -> 19558.3335
$ws.Cells.Item(89, 13).Value = -166673384   # M89: -277788234 -> -166673384
$ws.Cells.Item(89, 14).Value = -30790.3335   # N89: -31732 -> -30790.3335
$ws.Cells.Item(137, 8).Value = 3922155   # H137: 4651850 -> 3922155
$ws.Cells.Item(137, 9).Value = 498.52942   # I137: 619.5714 -> 498.52942
$ws.Cells.Item(137, 10).Value = 11765468   # J137: 13334146 -> 11765468
$ws.Cells.Item(137, 11).Value = 1495.58826   # K137: 1858.7142 -> 1495.58826
$ws.Cells.Item(137, 12).Value = 35296404   # L137: 40002438 -> 35296404
$ws.Cells.Item(137, 13).Value = 1054.41174   # M137: 691.2857999999999 -> 1054.41174
$ws.Cells.Item(137, 14).Value = -35301504   # N137: -40007538 -> -35301504

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 813.96295   # H2: 789.63635 -> 813.96295
$ws.Cells.Item(2, 9).Value = 884.1053000000001   # I2: 805.8889 -> 884.1053000000001
$ws.Cells.Item(2, 10).Value = 647.375   # J2: 716.5 -> 647.375
$ws.Cells.Item(2, 11).Value = 884.1053000000001   # K2: 805.8889 -> 884.1053000000001
$ws.Cells.Item(2, 12).Value = 647.375   # L2: 716.5 -> 647.375
$ws.Cells.Item(2, 13).Value = -771.1053000000001   # M2: -692.8889 -> -771.1053000000001
$ws.Cells.Item(2, 14).Value = -873.375   # N2: -942.5 -> -873.375
$ws.Cells.Item(61, 8).Value = 9435245   # H61: 10205472 -> 9435245
$ws.Cells.Item(61, 9).Value = 10205409   # I61: 11629378 -> 10205409
$ws.Cells.Item(61, 10).Value = 740.5   # J61: 810.3333 -> 740.5
$ws.Cells.Item(61, 11).Value = 10205409   # K61: 11629378 -> 10205409
$ws.Cells.Item(61, 12).Value = 740.5   # L61: 810.3333 -> 740.5
$ws.Cells.Item(61, 13).Value = -10205197   # M61: -11629166 -> -10205197
$ws.Cells.Item(61, 14).Value = -1164.5   # N61: -1234.3333 -> -1164.5
$ws.Cells.Item(74, 8).Value = 8198003   # H74: 9092324 -> 8198003
$ws.Cells.Item(74, 9).Value = 10417908   # I74: 11906162 -> 10417908
$ws.Cells.Item(74, 10).Value = 1429.5385   # J74: 1464.9231 -> 1429.5385
$ws.Cells.Item(74, 11).Value = 10417908   # K74: 11906162 -> 10417908
$ws.Cells.Item(74, 12).Value = 1429.5385   # L74: 1464.9231 -> 1429.5385
$ws.Cells.Item(74, 13).Value = -10417034   # M74: -11905288 -> -10417034
$ws.Cells.Item(74, 14).Value = -3177.5385   # N74: -3212.9231 -> -3177.5385
$ws.Cells.Item(77, 8).Value = 8198003   # H77: 9092324 -> 8198003
$ws.Cells.Item(77, 9).Value = 10417908   # I77: 11906162 -> 10417908
$ws.Cells.Item(77, 10).Value = 1429.5385   # J77: 1464.9231 -> 1429.5385
$ws.Cells.Item(77, 11).Value = 52089540   # K77: 59530810 -> 52089540
$ws.Cells.Item(77, 12).Value = 7147.692500000001   # L77: 7324.6155 -> 7147.692500000001
$ws.Cells.Item(77, 13).Value = -52085172   # M77: -59526442 -> -52085172
$ws.Cells.Item(77, 14).Value = -15883.6925   # N77: -16060.6155 -> -15883.6925
$ws.Cells.Item(116, 8).Value = 813.96295   # H116: 789.63635 -> 813.96295
$ws.Cells.Item(116, 9).Value = 884.1053000000001   # I116: 805.8889 -> 884.1053000000001
$ws.Cells.Item(116, 10).Value = 647.375   # J116: 716.5 -> 647.375
$ws.Cells.Item(116, 11).Value = 884.1053000000001   # K116: 805.8889 -> 884.1053000000001
$ws.Cells.Item(116, 12).Value = 647.375   # L116: 716.5 -> 647.375
$ws.Cells.Item(116, 13).Value = 1409.8947   # M116: 1488.1111 -> 1409.8947
$ws.Cells.Item(116, 14).Value = -5235.375   # N116: -5304.5 -> -5235.375
$ws.Cells.Item(136, 8).Value = 9435245   # H136: 10205472 -> 9435245
$ws.Cells.Item(136, 9).Value = 10205409   # I136: 11629378 -> 10205409
$ws.Cells.Item(136, 10).Value = 740.5   # J136: 810.3333 -> 740.5
$ws.Cells.Item(136, 11).Value = 30616227   # K136: 34888134 -> 30616227
$ws.Cells.Item(136, 12).Value = 2221.5   # L136: 2430.9999 -> 2221.5
$ws.Cells.Item(136, 13).Value = -30613677   # M136: -34885584 -> -30613677
$ws.Cells.Item(136, 14).Value = -7321.5   # N136: -7530.9999 -> -7321.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 813.96295   # H3: 789.63635 -> 813.96295
$ws.Cells.Item(3, 9).Value = 884.1053000000001   # I3: 805.8889 -> 884.1053000000001
$ws.Cells.Item(3, 10).Value = 647.375   # J3: 716.5 -> 647.375
$ws.Cells.Item(3, 11).Value = 884.1053000000001   # K3: 805.8889 -> 884.1053000000001
$ws.Cells.Item(3, 12).Value = 647.375   # L3: 716.5 -> 647.375
$ws.Cells.Item(3, 13).Value = -770.1053000000001   # M3: -691.8889 -> -770.1053000000001
$ws.Cells.Item(3, 14).Value = -875.375   # N3: -944.5 -> -875.375
$ws.Cells.Item(63, 8).Value = 48522.6   # H63: 47647.332 -> 48522.6
$ws.Cells.Item(63, 10).Value = 48522.6   # J63: 47647.332 -> 48522.6
$ws.Cells.Item(63, 12).Value = 48522.6   # L63: 47647.332 -> 48522.6
$ws.Cells.Item(63, 14).Value = -49894.6   # N63: -49019.332 -> -49894.6
$ws.Cells.Item(66, 8).Value = 48522.6   # H66: 47647.332 -> 48522.6
$ws.Cells.Item(66, 10).Value = 48522.6   # J66: 47647.332 -> 48522.6
$ws.Cells.Item(66, 12).Value = 145567.8   # L66: 142941.996 -> 145567.8
$ws.Cells.Item(66, 14).Value = -152431.8   # N66: -149805.996 -> -152431.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(20, 8).Value = 33150   # H20: 35450 -> 33150
$ws.Cells.Item(20, 10).Value = 33150   # J20: 35450 -> 33150
$ws.Cells.Item(20, 12).Value = 33150   # L20: 35450 -> 33150
$ws.Cells.Item(20, 14).Value = -33622   # N20: -35922 -> -33622
$ws.Cells.Item(30, 8).Value = 33150   # H30: 35450 -> 33150
$ws.Cells.Item(30, 10).Value = 33150   # J30: 35450 -> 33150
$ws.Cells.Item(30, 12).Value = 33150   # L30: 35450 -> 33150
$ws.Cells.Item(30, 14).Value = -33332   # N30: -35632 -> -33332
$ws.Cells.Item(31, 8).Value = 6537938.5   # H31: 6668667 -> 6537938.5
$ws.Cells.Item(31, 9).Value = 1819.7234   # I31: 1808.8937 -> 1819.7234
$ws.Cells.Item(31, 10).Value = 83337336   # J31: 111116110 -> 83337336
$ws.Cells.Item(31, 11).Value = 1819.7234   # K31: 1808.8937 -> 1819.7234
$ws.Cells.Item(31, 12).Value = 83337336   # L31: 111116110 -> 83337336
$ws.Cells.Item(31, 13).Value = -1524.7234   # M31: -1513.8937 -> -1524.7234
$ws.Cells.Item(31, 14).Value = -83337926   # N31: -111116700 -> -83337926
$ws.Cells.Item(34, 8).Value = 6537938.5   # H34: 6668667 -> 6537938.5
$ws.Cells.Item(34, 9).Value = 1819.7234   # I34: 1808.8937 -> 1819.7234
$ws.Cells.Item(34, 10).Value = 83337336   # J34: 111116110 -> 83337336
$ws.Cells.Item(34, 11).Value = 1819.7234   # K34: 1808.8937 -> 1819.7234
$ws.Cells.Item(34, 12).Value = 83337336   # L34: 111116110 -> 83337336
$ws.Cells.Item(34, 13).Value = -1617.7234   # M34: -1606.8937 -> -1617.7234
$ws.Cells.Item(34, 14).Value = -83337740   # N34: -111116514 -> -83337740
$ws.Cells.Item(128, 8).Value = 33150   # H128: 35450 -> 33150
$ws.Cells.Item(128, 10).Value = 33150   # J128: 35450 -> 33150
$ws.Cells.Item(128, 12).Value = 33150   # L128: 35450 -> 33150
$ws.Cells.Item(128, 14).Value = -43110   # N128: -45410 -> -43110
$ws.Cells.Item(134, 8).Value = 1359.1094   # H134: 1480.2222 -> 1359.1094
$ws.Cells.Item(134, 9).Value = 1270.1187   # I134: 1375.74 -> 1270.1187
$ws.Cells.Item(134, 10).Value = 2409.2   # J134: 2786.25 -> 2409.2
$ws.Cells.Item(134, 11).Value = 3810.3561   # K134: 4127.22 -> 3810.3561
$ws.Cells.Item(134, 12).Value = 7227.599999999999   # L134: 8358.75 -> 7227.599999999999
$ws.Cells.Item(134, 13).Value = -1275.3561   # M134: -1592.22 -> -1275.3561
$ws.Cells.Item(134, 14).Value = -12297.6   # N134: -13428.75 -> -12297.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 902.94116   # H5: 778.8095 -> 902.94116
$ws.Cells.Item(5, 9).Value = 392.375   # I5: 345.33334 -> 392.375
$ws.Cells.Item(5, 11).Value = 1177.125   # K5: 1036.00002 -> 1177.125
$ws.Cells.Item(5, 13).Value = -1065.125   # M5: -924.0000199999999 -> -1065.125
$ws.Cells.Item(132, 8).Value = 3783.3333   # H132: 1441.9166 -> 3783.3333
$ws.Cells.Item(132, 9).Value = 500   # I132: 436 -> 500
$ws.Cells.Item(132, 10).Value = 5425   # J132: 1944.875 -> 5425
$ws.Cells.Item(132, 11).Value = 4500   # K132: 3924 -> 4500
$ws.Cells.Item(132, 12).Value = 48825   # L132: 17503.875 -> 48825
$ws.Cells.Item(132, 13).Value = -1970   # M132: -1394 -> -1970
$ws.Cells.Item(132, 14).Value = -53885   # N132: -22563.875 -> -53885
$ws.Cells.Item(135, 8).Value = 902.94116   # H135: 778.8095 -> 902.94116
$ws.Cells.Item(135, 9).Value = 392.375   # I135: 345.33334 -> 392.375
$ws.Cells.Item(135, 11).Value = 3531.375   # K135: 3108.00006 -> 3531.375
$ws.Cells.Item(135, 13).Value = -996.375   # M135: -573.0000600000003 -> -996.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 4180.154   # H132: 2953.1772 -> 4180.154
$ws.Cells.Item(132, 9).Value = 3446.7   # I132: 2159.709 -> 3446.7
$ws.Cells.Item(132, 10).Value = 5180.3184   # J132: 4771.5415 -> 5180.3184
$ws.Cells.Item(132, 11).Value = 10340.1   # K132: 6479.126999999999 -> 10340.1
$ws.Cells.Item(132, 12).Value = 15540.9552   # L132: 14314.6245 -> 15540.9552
$ws.Cells.Item(132, 13).Value = -7810.099999999999   # M132: -3949.126999999999 -> -7810.099999999999
$ws.Cells.Item(132, 14).Value = -20600.9552   # N132: -19374.6245 -> -20600.9552

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 6762383   # H132: 7358933.5 -> 6762383
$ws.Cells.Item(132, 9).Value = 3033.5254   # I132: 3408.9424 -> 3033.5254
$ws.Cells.Item(132, 10).Value = 33349158   # J132: 31264388 -> 33349158
$ws.Cells.Item(132, 11).Value = 9100.5762   # K132: 10226.8272 -> 9100.5762
$ws.Cells.Item(132, 12).Value = 100047474   # L132: 93793164 -> 100047474
$ws.Cells.Item(132, 13).Value = -6570.5762   # M132: -7696.8272 -> -6570.5762
$ws.Cells.Item(132, 14).Value = -100052534   # N132: -93798224 -> -100052534

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 19740.4   # H62: 11014.5 -> 19740.4
$ws.Cells.Item(62, 9).Value = 6933   # I62: 7000 -> 6933
$ws.Cells.Item(62, 10).Value = 38951.5   # J62: 12109.363 -> 38951.5
$ws.Cells.Item(62, 11).Value = 6933   # K62: 7000 -> 6933
$ws.Cells.Item(62, 12).Value = 38951.5   # L62: 12109.363 -> 38951.5
$ws.Cells.Item(62, 13).Value = -6309   # M62: -6376 -> -6309
$ws.Cells.Item(62, 14).Value = -40199.5   # N62: -13357.363 -> -40199.5
$ws.Cells.Item(65, 8).Value = 19740.4   # H65: 11014.5 -> 19740.4
$ws.Cells.Item(65, 9).Value = 6933   # I65: 7000 -> 6933
$ws.Cells.Item(65, 10).Value = 38951.5   # J65: 12109.363 -> 38951.5
$ws.Cells.Item(65, 11).Value = 34665   # K65: 35000 -> 34665
$ws.Cells.Item(65, 12).Value = 194757.5   # L65: 60546.815 -> 194757.5
$ws.Cells.Item(65, 13).Value = -31545   # M65: -31880 -> -31545
$ws.Cells.Item(65, 14).Value = -200997.5   # N65: -66786.815 -> -200997.5
$ws.Cells.Item(113, 8).Value = 1379.5   # H113: 1516.1904 -> 1379.5
$ws.Cells.Item(113, 9).Value = 1210.25   # I113: 1154.5834 -> 1210.25
$ws.Cells.Item(113, 10).Value = 1548.75   # J113: 1998.3334 -> 1548.75
$ws.Cells.Item(113, 11).Value = 3630.75   # K113: 3463.7502 -> 3630.75
$ws.Cells.Item(113, 12).Value = 4646.25   # L113: 5995.0002 -> 4646.25
$ws.Cells.Item(113, 13).Value = -1460.75   # M113: -1293.7502 -> -1460.75
$ws.Cells.Item(113, 14).Value = -8986.25   # N113: -10335.0002 -> -8986.25
